$p = $ppt.ActivePresentation
$newDate = "08/02/2022"
$ppPlaceholderDate = 16

function Set-DatePlaceholderText($shapes, $text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $isDatePh = $false
            if ($sh.PlaceholderFormat -ne $null) {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            }
            if ($isDatePh) {
                $sh.TextFrame.TextRange.Text = $text
            }
        }
    }
}

# Update the slide master's date placeholder.
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

# Update every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}
